# Fix js calendar. Fix Bozena and 'Name, Surname' in resources
#
# 1) Normalize "Giermek Bozena" (no diacritic) -> "Giermek Bozena" with the
#    proper Polish "z with dot above" (z-with-dot, U+017C) so it matches the
#    already-existing "Giermek Bo\u017Cena" professor entry.
# 2) Swap "Iwaniec Joanna" (Surname Name) -> "Joanna Iwaniec" (Name Surname)
#    for every occurrence in the Professor column.
# 3) Restore workbook view/selection state: Sheet1 scrolled down with E78
#    selected, Sheet2/Sheet3 with A1 selected (instead of whole column A),
#    and Sheet1 left as the active/selected tab.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

$wrongBozena = "Giermek Bozena"
$correctBozena = "Giermek Bo" + [char]0x017C + "ena"
$wrongIwaniec = "Iwaniec Joanna"
$correctIwaniec = "Joanna Iwaniec"

$lastRow = $ws1.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws1.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq $wrongBozena) {
        $cell.Value2 = $correctBozena
    } elseif ($val -eq $wrongIwaniec) {
        $cell.Value2 = $correctIwaniec
    }
}

$win = $excel.ActiveWindow

# Sheet2: select A1 (was whole column A), keep it not the active tab.
$ws2.Activate() | Out-Null
$win.DisplayGridlines = $true
$ws2.Range("A1").Select() | Out-Null

# Sheet3: select A1 (was whole column A), keep it not the active tab.
$ws3.Activate() | Out-Null
$win.DisplayGridlines = $true
$ws3.Range("A1").Select() | Out-Null

# Sheet1: becomes/remains the active tab, scrolled to row 51, E78 selected.
$ws1.Activate() | Out-Null
$win.DisplayGridlines = $true
$win.ScrollRow = 51
$win.ScrollColumn = 1
$ws1.Range("E78").Select() | Out-Null
